# "added 64 bit designs" - refresh the utilization report numbers for the
# only_integer64 / 30mhz / mxu_5x5 configuration, and widen column F to match
# the other data columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (6th column) widens from the narrow (10.46875) preset to the wider
# (11.71875) preset used by columns B, D and E. The host's ColumnWidth setter
# quantizes to 1/6-character steps, so feed it the mid-bucket value that
# resolves to the closest achievable width (11.666666666666666 ~= 11.71875).
$ws.Columns.Item(6).ColumnWidth = 10.833333333333334

# Updated utilization percentages (row 2) for the new 64-bit designs.
$ws.Range("B2").Value = 82.06015014648438
$ws.Range("C2").Value = 6.0804595947265625
$ws.Range("D2").Value = 24.877819061279297
$ws.Range("E2").Value = 57.85714340209961
$ws.Range("F2").Value = 95.45454406738281
